# Tracker auto-update: append one new result row (row 12) to Sheet1,
# matching the new A1:H12 dimension.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 14655434

# Force the date-looking value to be stored as literal text (not an
# Excel date serial) by switching the cell to a text number format
# before assigning it, then dropping back to the default style so no
# visible formatting change is left behind on the cell.
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2025-09-18"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").Value = "Beatriz Haddad Maia"
$ws.Range("D12").Value = "Ella Seidel"
$ws.Range("E12").Value = "Gana Beatriz Haddad Maia"
$ws.Range("F12").Value = 1.73

# resultado / profit are still blank for this not-yet-played match.
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
